# Add "biomass" to the least-cost-dispatch electricity source subscript,
# and flag on the About sheet that nuclear is intentionally excluded from
# dispatch (it's handled as guaranteed dispatch elsewhere).

$wb = $excel.ActiveWorkbook

# --- About sheet: explain why "nuclear es" isn't in the dispatch list ---
$wsAbout = $wb.Worksheets.Item("About")

# Highlight "nuclear es" (A13) in red.
$wsAbout.Range("A13").Font.Color = 255

# Add the explanatory note below the existing source list (row 34 stays blank).
$wsAbout.Range("A35").Value = "We do not use nuclear for dispatch in the US because we have nuclear listed"
$wsAbout.Range("A36").Value = "as guaranteed dispatch in elec/BGDPbES."

# --- ESUfRaLCD-dispatch sheet: insert "biomass" as a dispatch source ---
$wsDispatch = $wb.Worksheets.Item("ESUfRaLCD-dispatch")

# Insert a new row above the old row 5 ("petroleum"), shifting everything
# below it down by one and extending the shared "... es"/"... dispatch"
# formulas automatically.
$wsDispatch.Rows("5:5").Insert()

# Fill in the new biomass row with literal values (matching how "nuclear es"
# is handled elsewhere in this workbook - a plain value, not a formula).
$wsDispatch.Range("A5").Value = "biomass"
$wsDispatch.Range("B5").Value = "biomass es"
$wsDispatch.Range("C5").Value = "biomass dispatch"

# --- Restore the selection/active-sheet state ---
$wsDispatch.Activate() | Out-Null
$wsDispatch.Range("A15").Select() | Out-Null

$wsAbout.Activate() | Out-Null
$wsAbout.Range("C26").Select() | Out-Null
